$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell $ws 'D2' '64.287.41'
Set-TextCell $ws 'E2' '  -0.72%  '
Set-TextCell $ws 'D3' '3.345.70'
Set-TextCell $ws 'E3' '  -2.61%  '
Set-TextCell $ws 'D5' '555.17'
Set-TextCell $ws 'E5' '  -2.64%  '
Set-TextCell $ws 'D6' '175.03'
Set-TextCell $ws 'E6' '  +0.08%  '
Set-TextCell $ws 'D7' '0.618'
Set-TextCell $ws 'E7' '  -0.76%  '
Set-TextCell $ws 'D8' '3.334.94'
Set-TextCell $ws 'E8' '  -2.70%  '
Set-TextCell $ws 'E9' '  +0.03%  '
Set-TextCell $ws 'D10' '0.627'
Set-TextCell $ws 'E10' '  +0.92%  '
Set-TextCell $ws 'D11' '0.162'
Set-TextCell $ws 'E11' '  +2.59%  '
Set-TextCell $ws 'D12' '54.43'
Set-TextCell $ws 'E12' '  -0.55%  '
Set-TextCell $ws 'D13' '0.0000272'
Set-TextCell $ws 'E13' '  +0.28%  '
Set-TextCell $ws 'D14' '9.06'
Set-TextCell $ws 'E14' '  -0.30%  '
Set-TextCell $ws 'D15' '3.893.81'
Set-TextCell $ws 'E15' '  -2.32%  '
Set-TextCell $ws 'D16' '18.38'
Set-TextCell $ws 'E16' '  +1.97%  '
Set-TextCell $ws 'E17' '  -2.07%  '
Set-TextCell $ws 'D18' '3.360.73'
Set-TextCell $ws 'E18' '  -2.41%  '
Set-TextCell $ws 'D19' '11.80'
Set-TextCell $ws 'E19' '  -0.13%  '
Set-TextCell $ws 'D20' '64.336.18'
Set-TextCell $ws 'E20' '  -0.72%  '
Set-TextCell $ws 'D21' '0.981'
Set-TextCell $ws 'E21' '  -0.29%  '
Set-TextCell $ws 'D22' '454.07'
Set-TextCell $ws 'E22' '  +11.77%  '
Set-TextCell $ws 'D23' '4.88'
Set-TextCell $ws 'E23' '  +11.43%  '
Set-TextCell $ws 'D24' '4.07'
Set-TextCell $ws 'E24' '  -2.28%  '
Set-TextCell $ws 'D25' '85.65'
Set-TextCell $ws 'E25' '  +2.60%  '
Set-TextCell $ws 'E26' '  -0.67%  '
Set-TextCell $ws 'D27' '10.92'
Set-TextCell $ws 'E27' '  +1.47%  '
Set-TextCell $ws 'D28' '2.84'
Set-TextCell $ws 'E28' '  +1.85%  '
Set-TextCell $ws 'D29' '8.75'
Set-TextCell $ws 'E29' '  -2.12%  '
Set-TextCell $ws 'D30' '29.92'
Set-TextCell $ws 'E30' '  +0.32%  '
Set-TextCell $ws 'E31' '  +0.71%  '
Set-TextCell $ws 'B32' 'Bittensor'
Set-TextCell $ws 'C32' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextCell $ws 'D32' '584.13'
Set-TextCell $ws 'E32' '  +0.42%  '
Set-TextCell $ws 'B33' 'Cosmos'
Set-TextCell $ws 'C33' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell $ws 'D33' '11.46'
Set-TextCell $ws 'E33' '  -0.38%  '
Set-TextCell $ws 'E34' '  -0.32%  '
Set-TextCell $ws 'D35' '58.54'
Set-TextCell $ws 'E35' '  -1.62%  '
Set-TextCell $ws 'D36' '0.999'
Set-TextCell $ws 'E36' '  -0.02%  '
Set-TextCell $ws 'D37' '0.140'
Set-TextCell $ws 'E37' '  -7.78%  '
Set-TextCell $ws 'B38' 'Stacks'
Set-TextCell $ws 'C38' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextCell $ws 'D38' '3.49'
Set-TextCell $ws 'E38' '  -1.60%  '
Set-TextCell $ws 'B39' 'InjectiveProtocol'
Set-TextCell $ws 'C39' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextCell $ws 'D39' '35.68'
Set-TextCell $ws 'E39' '  -1.18%  '
Set-TextCell $ws 'D40' '0.0₃0754'
Set-TextCell $ws 'E40' '  -2.20%  '
Set-TextCell $ws 'D41' '0.374'
Set-TextCell $ws 'E41' '  -0.14%  '
Set-TextCell $ws 'D42' '3.095.58'
Set-TextCell $ws 'E42' '  -2.18%  '
Set-TextCell $ws 'D43' '1.00'
Set-TextCell $ws 'E43' '  +0.15%  '
Set-TextCell $ws 'B44' 'ThetaToken'
Set-TextCell $ws 'C44' 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextCell $ws 'D44' '2.79'
Set-TextCell $ws 'E44' '  -4.55%  '
Set-TextCell $ws 'B45' 'Fetch.AI'
Set-TextCell $ws 'C45' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextCell $ws 'D45' '2.52'
Set-TextCell $ws 'E45' '  +1.07%  '
Set-TextCell $ws 'D46' '3.22'
Set-TextCell $ws 'E46' '  -0.90%  '
Set-TextCell $ws 'D47' '0.0409'
Set-TextCell $ws 'E47' '  +0.26%  '
Set-TextCell $ws 'D48' '0.130'
Set-TextCell $ws 'E48' '  +0.03%  '
Set-TextCell $ws 'D49' '2.57'
Set-TextCell $ws 'E49' '  -2.21%  '
Set-TextCell $ws 'B50' 'THORChain'
Set-TextCell $ws 'C50' 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextCell $ws 'D50' '8.31'
Set-TextCell $ws 'E50' '  -1.47%  '
Set-TextCell $ws 'B51' 'Monero'
Set-TextCell $ws 'C51' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell $ws 'D51' '135.54'
Set-TextCell $ws 'E51' '  -0.98%  '
